$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge any existing merged cells on the sheet before rewriting content
foreach ($mc in @("A2:C2","A5:C5","A8:C8","A11:C11","A16:C16","A20:C20","A22:C22","A25:C25")) {
    $ws.Range($mc).UnMerge()
}

# Clear the previously used range entirely (values + formatting leftovers)
$ws.Range("A1:C26").Clear()

# Header row
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "Quantity"
$ws.Range("C1").Value = "Revenue ▼"

# Lassi section
$ws.Range("A2").Value = "Lassi"
$ws.Range("A3").Value = "Mango Lassi"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "₹480.00"
$ws.Range("A4").Value = "Butterscotch Lassi"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "₹60.00"
$ws.Range("A5").Value = "Strawberry Lassi"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "₹40.00"

# Burgers section
$ws.Range("A6").Value = "Burgers"
$ws.Range("A7").Value = "Chicken Cheese Burger"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "₹180.00"
$ws.Range("A8").Value = "Chicken Burger"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "₹120.00"

# Pops section
$ws.Range("A9").Value = "Pops"
$ws.Range("A10").Value = "Veg Cheese Pops"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "₹70.00"
$ws.Range("A11").Value = "Chicken Cheese Pops"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "₹70.00"

# Wraps section
$ws.Range("A12").Value = "Wraps"
$ws.Range("A13").Value = "Chicken Wrap"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "₹60.00"

# Milkshakes section
$ws.Range("A14").Value = "Milkshakes"
$ws.Range("A15").Value = "Oreo Shake"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "₹50.00"
$ws.Range("A16").Value = "Banana Shake"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "₹50.00"
$ws.Range("A17").Value = "Vanilla Shake"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "₹40.00"

# Re-merge the category header rows
$ws.Range("A2:C2").Merge()
$ws.Range("A6:C6").Merge()
$ws.Range("A9:C9").Merge()
$ws.Range("A12:C12").Merge()
$ws.Range("A14:C14").Merge()
